$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values in row 2 as per the FlashScore data refresh
$ws.Range("G2").Value = 2.35    # Odd_H_FT: 2.3 -> 2.35
$ws.Range("I2").Value = 2.88    # Odd_A_FT: 2.9 -> 2.88
$ws.Range("L2").Value = 3.2     # Odd_A_HT: 3.25 -> 3.2
$ws.Range("N2").Value = 17      # Odd_Under05_FT: 15 -> 17
$ws.Range("X2").Value = 15      # Odd_CS_2-0: 13 -> 15
$ws.Range("AP2").Value = 19     # Odd_CS_2-1_HT: 17 -> 19
$ws.Range("AW2").Value = 351    # Odd_CS_3-3_HT: 301 -> 351
